# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Estado de Cuenta" detail table (rows 16-29, columns B:G) is refreshed
# with new worker / period data. Two workers, each with 7 overdue periods
# (2010, 2011, 2012, 2101, 2102, 2103, 2104), now sorted with the most
# recent period (2104) first, "URBIS AMAYA ACONCHA" (CC 8643805) block
# followed by "LUIS FERNANDO VANEGAS BASANTA" (CC 1042427825) block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("CC", "8643805",    "URBIS AMAYA ACONCHA",           "2104", 59322, 1711200),
    @("CC", "8643805",    "URBIS AMAYA ACONCHA",           "2103", 68448, 1711200),
    @("CC", "8643805",    "URBIS AMAYA ACONCHA",           "2102", 68448, 1711200),
    @("CC", "8643805",    "URBIS AMAYA ACONCHA",           "2101", 68448, 1711200),
    @("CC", "8643805",    "URBIS AMAYA ACONCHA",           "2012", 68448, 1711200),
    @("CC", "8643805",    "URBIS AMAYA ACONCHA",           "2011", 68448, 1711200),
    @("CC", "8643805",    "URBIS AMAYA ACONCHA",           "2010", 68448, 1711200),
    @("CC", "1042427825", "LUIS FERNANDO VANEGAS BASANTA", "2104", 34965, 1008600),
    @("CC", "1042427825", "LUIS FERNANDO VANEGAS BASANTA", "2103", 40344, 1008600),
    @("CC", "1042427825", "LUIS FERNANDO VANEGAS BASANTA", "2102", 40344, 1008600),
    @("CC", "1042427825", "LUIS FERNANDO VANEGAS BASANTA", "2101", 40344, 1008600),
    @("CC", "1042427825", "LUIS FERNANDO VANEGAS BASANTA", "2012", 40344, 1008600),
    @("CC", "1042427825", "LUIS FERNANDO VANEGAS BASANTA", "2011", 40344, 1008600),
    @("CC", "1042427825", "LUIS FERNANDO VANEGAS BASANTA", "2010", 40344, 1008600)
)

$r = 16
foreach ($row in $data) {
    $ws.Range("B$r").Value = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $ws.Range("G$r").Value = $row[5]
    $r = $r + 1
}
